$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 11, pushing "register" and "counter" rows down to 12 and 13.
$ws.Rows("11:11").Insert()

# --- Row 10 (rom) : add part number / note / datasheet link ---
$ws.Range("B10").Value = "MT29F8G08ABACA"
$ws.Range("C10").Value = "same addr and i/o"
$ws.Hyperlinks.Add($ws.Range("D10"), "https://datasheet.lcsc.com/lcsc/1912111437_Micron-Tech-MT29F8G08ABACAWP-IT-C_C400999.pdf") | Out-Null
$ws.Range("D10").Style = "Hyperlink"
$ws.Range("D10").HorizontalAlignment = 5

# --- Row 11 (new, blank Function) : AT28C64 ---
$ws.Range("B11").Value = "AT28C64"
$ws.Range("D11").Value = "https://ww1.microchip.com/downloads/en/DeviceDoc/doc0270.pdf"
$ws.Range("D11").Style = "Hyperlink"
$ws.Range("D11").HorizontalAlignment = 5

# --- Row 12 (register, shifted down from 11) : SN74HCT574PWR ---
$url12 = "https://www.ti.com/lit/ds/symlink/sn74hct574.pdf?ts=1640527510946&ref_url=https%253A%252F%252Fwww.ti.com%252Fproduct%252FSN74HCT574%253Futm_source%253Dgoogle%2526utm_medium%253Dcpc%2526utm_campaign%253Dasc-int-lvt-prodfolderdynamic-cpc-pf-google-wwe%2526utm_content%253Dprodfolddynamic%2526ds_k%253DDYNAMIC%2BSEARCH%2BADS%2526DCM%253Dyes%2526gclid%253DCj0KCQiAwqCOBhCdARIsAEPyW9k3RiLEWJBSyZsyJfcgBFzijBS9H81wBqR6zBzhKHt4lWJ7ZB4cynsaAolWEALw_wcB%2526gclsrc%253Daw.ds"
$ws.Range("B12").Value = "SN74HCT574PWR"
$ws.Hyperlinks.Add($ws.Range("D12"), $url12, "", "", $url12) | Out-Null
$ws.Range("D12").Style = "Hyperlink"

# --- Row 13 (counter, shifted down from 12) : SN74HC161DR ---
$ws.Range("A13").Value = "counter"
$ws.Range("B13").Value = "SN74HC161DR"
$ws.Range("C13").Value = "4 bit only"
$ws.Hyperlinks.Add($ws.Range("D13"), "https://www.ti.com/lit/ds/symlink/sn74hc161.pdf?HQS=dis-mous-null-mousermode-dsf-pf-null-wwe&ts=1640512735338") | Out-Null
$ws.Range("D13").Style = "Hyperlink"
$ws.Range("D13").HorizontalAlignment = 5

# Match the saved selection state from the edited workbook.
$ws.Range("C14:C15").Select()
